# Resize the worksheet's columns:
#   - Columns A:C -> uniform width (~22.71 "characters", i.e. stored OOXML width 22.7109375)
#   - Column D   -> width grows by one character step (~47.71, stored OOXML width 47.7109375)
#   - Column E   -> width grows by one character step (~22.71, stored OOXML width 22.7109375)
#
# Note: this runtime's ColumnWidth setter quantizes to its own internal pixel
# grid, so the literal "characters" values below are chosen so that, after
# that quantization, the saved <col .../> width lands as close as possible
# to the intended 22.7109375 / 47.7109375 target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetWidthNarrow = 21.833333333333332   # -> saved width ~22.7109375
$targetWidthWide   = 46.833333333333336   # -> saved width ~47.7109375

$ws.Range("A1:C1").EntireColumn.ColumnWidth = $targetWidthNarrow
$ws.Range("D1").EntireColumn.ColumnWidth = $targetWidthWide
$ws.Range("E1").EntireColumn.ColumnWidth = $targetWidthNarrow
